$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 18 (CreateDate): change 形態 (Type) from DATE to TIMESTAMP and clear 長度 (Length)
$ws.Range("D18").Value = "TIMESTAMP"
$ws.Range("E18").Value = ""

# Row 20 (LastUpdate): change 形態 (Type) from DATE to TIMESTAMP and clear 長度 (Length)
$ws.Range("D20").Value = "TIMESTAMP"
$ws.Range("E20").Value = ""

# Leave the selection on the last-edited cell, matching the saved view state
$ws.Activate()
$ws.Range("E20").Select()
